# Update prediction data: remove the stale "Ezra Mayers" row.
# The rows below it (Alphonse Aréola, Mads Hermansen) shift up to take
# its place, and the sheet's used range shrinks by one row
# (A1:DL23 -> A1:DL22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).Delete()
